# Add "PerItem" per-key columns (R:U) + medians (row 30) + speedup ratios (row 31)
# to both ReadTests and WriteTests sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ReadTests")
$ws2 = $wb.Worksheets.Item("WriteTests")

# ---- Headers (row 23) ----
$ws1.Range("R23").Value = "PerItem Simple"
$ws1.Range("S23").Value = "PerItem Expected"
$ws1.Range("T23").Value = "PerItem Parallel"
$ws1.Range("U23").Value = "PerItem Bulk"

$ws2.Range("R23").Value = "PerItem Simple"
$ws2.Range("S23").Value = "PerItem Expected"
$ws2.Range("T23").Value = "PerItem Parallel"
$ws2.Range("U23").Value = "PerItem Bulk"

# ---- ReadTests (sheet1): per-item columns, rows 24-29 ----
$ws1.Range("R24").Formula = "=J24/`$I24"
$ws1.Range("S24").Formula = "=K24/`$I24"
$ws1.Range("T24").Formula = "=L24/`$I24"
$ws1.Range("U24").Formula = "=M24/`$I24"

for ($r = 25; $r -le 29; $r++) {
    $ws1.Range("R$r").Formula = "=J$r/`$I$r"
    $ws1.Range("S$r").Formula = "=K$r/`$I$r"
    $ws1.Range("T$r").Formula = "=L$r/`$I$r"
    $ws1.Range("U$r").Formula = "=M$r/`$I$r"
}

# Row 30: medians
$ws1.Range("R30").Formula = "=MEDIAN(R24:R29)"
$ws1.Range("S30").Formula = "=MEDIAN(S24:S29)"
$ws1.Range("T30").Formula = "=MEDIAN(T24:T29)"
$ws1.Range("U30").Formula = "=MEDIAN(U24:U29)"

# Row 31: speedup ratios
$ws1.Range("T31").Formula = "=R30/T30"
$ws1.Range("U31").Formula = "=S30/U30"

# ---- WriteTests (sheet2): per-item columns, rows 24-29 ----
$ws2.Range("R24").Formula = "=J24/`$I24"
$ws2.Range("S24").Formula = "=K24/`$I24"
$ws2.Range("T24").Formula = "=L24/`$I24"
$ws2.Range("U24").Formula = "=M24/`$I24"

# Note: row 25's R formula uses a relative I25 (no $), matching the source workbook.
$ws2.Range("R25").Formula = "=J25/I25"
for ($r = 26; $r -le 29; $r++) {
    $ws2.Range("R$r").Formula = "=J$r/I$r"
}
for ($r = 25; $r -le 29; $r++) {
    $ws2.Range("S$r").Formula = "=K$r/`$I$r"
    $ws2.Range("T$r").Formula = "=L$r/`$I$r"
    $ws2.Range("U$r").Formula = "=M$r/`$I$r"
}

# Row 30: medians
$ws2.Range("R30").Formula = "=MEDIAN(R24:R29)"
$ws2.Range("S30").Formula = "=MEDIAN(S24:S29)"
$ws2.Range("T30").Formula = "=MEDIAN(T24:T29)"
$ws2.Range("U30").Formula = "=MEDIAN(U24:U29)"

# Row 31: speedup ratios (absolute reference to the PerItem Simple median)
$ws2.Range("T31").Formula = "=`$R`$30/T30"
$ws2.Range("U31").Formula = "=`$R`$30/U30"

# WriteTests new cells carry a "0.00" number format (only the cells that
# actually hold a formula -- avoid materialising empty R31/S31 cells).
$ws2.Range("R24:U30").NumberFormat = "0.00"
$ws2.Range("T31:U31").NumberFormat = "0.00"

# ---- Column widths (approximate bestFit look for the new columns) ----
$ws1.Range("R1").ColumnWidth = 14.02
$ws1.Range("S1").ColumnWidth = 16.02
$ws1.Range("T1").ColumnWidth = 14.6
$ws1.Range("U1").ColumnWidth = 11.6

$ws2.Range("R1").ColumnWidth = 14.02
$ws2.Range("S1").ColumnWidth = 16.02
$ws2.Range("T1").ColumnWidth = 14.6
$ws2.Range("U1").ColumnWidth = 11.6

# ---- Selection / active sheet ----
# WriteTests keeps a selection but ReadTests becomes the active tab (matches the
# tabSelected flag moving from WriteTests to ReadTests in the saved file).
$ws2.Range("R31").Select()
$ws1.Activate()
$ws1.Range("R32").Select()

Write-Host "Applied PerItem columns + medians + speedup ratios"
